# Updated cryptos list (price & 1h volume refresh, plus a couple of
# ranking swaps) as scraped on Sat Jun  3 02:17:36 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link swaps (ranking reshuffled) --------------------
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'

# --- Price column updates ---------------------------------------------
# Some prices are plain numeric-looking strings (e.g. "0.5220") that must
# stay TEXT (leading/trailing zeros matter) - force Text format first,
# then drop back to the Normal style so no stray style index is left
# behind on the cell.
$priceTextCells = @(
    "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D17", "D18", "D19", "D21",
    "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D37",
    "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D48", "D49", "D50", "D51"
)
foreach ($addr in $priceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = '306.84'
$ws.Range("D7").Value = '0.5220'
$ws.Range("D8").Value = '0.3767'
$ws.Range("D9").Value = '0.07242'
$ws.Range("D10").Value = '21.21'
$ws.Range("D11").Value = '0.8983'
$ws.Range("D12").Value = '0.08340'
$ws.Range("D14").Value = '94.60'
$ws.Range("D15").Value = '5.273'
$ws.Range("D17").Value = '0.000008595'
$ws.Range("D18").Value = '14.51'
$ws.Range("D19").Value = '1.002'
$ws.Range("D21").Value = '5.058'
$ws.Range("D23").Value = '10.59'
$ws.Range("D24").Value = '6.427'
$ws.Range("D25").Value = '146.70'
$ws.Range("D26").Value = '2.282'
$ws.Range("D27").Value = '1.753'
$ws.Range("D28").Value = '18.15'
$ws.Range("D29").Value = '114.83'
$ws.Range("D30").Value = '4.940'
$ws.Range("D31").Value = '4.790'
$ws.Range("D32").Value = '0.09210'
$ws.Range("D33").Value = '0.8104'
$ws.Range("D35").Value = '1.239'
$ws.Range("D37").Value = '3.358'
$ws.Range("D38").Value = '2.560'
$ws.Range("D39").Value = '0.5691'
$ws.Range("D40").Value = '0.01972'
$ws.Range("D41").Value = '1.072'
$ws.Range("D42").Value = '8.969'
$ws.Range("D43").Value = '6.593'
$ws.Range("D44").Value = '118.36'
$ws.Range("D46").Value = '0.4820'
$ws.Range("D48").Value = '10.12'
$ws.Range("D49").Value = '1.610'
$ws.Range("D50").Value = '37.51'
$ws.Range("D51").Value = '63.53'

foreach ($addr in $priceTextCells) {
    $ws.Range($addr).Style = "Normal"
}

# Prices that are not ambiguous (contain multiple "." separators) can be
# assigned directly - Excel cannot coerce them to a number.
$ws.Range("D2").Value = '27.195.38'
$ws.Range("D3").Value = '1.905.03'
$ws.Range("D13").Value = '1.904.16'
$ws.Range("D20").Value = '27.243.56'
$ws.Range("D22").Value = '2.145.45'

# --- Volume(1h) column updates -----------------------------------------
# These always carry leading/trailing padding spaces so Excel keeps them
# as text automatically.
$ws.Range("E2").Value = '  +1.32%  '
$ws.Range("E3").Value = '  +1.83%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("E7").Value = '  +2.06%  '
$ws.Range("E8").Value = '  +2.43%  '
$ws.Range("E9").Value = '  +0.94%  '
$ws.Range("E10").Value = '  +2.81%  '
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("E12").Value = '  +10.73%  '
$ws.Range("E13").Value = '  +1.80%  '
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("E15").Value = '  +0.68%  '
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("E17").Value = '  +2.42%  '
$ws.Range("E18").Value = '  +2.05%  '
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("E20").Value = '  +1.35%  '
$ws.Range("E21").Value = '  +0.87%  '
$ws.Range("E22").Value = '  +1.94%  '
$ws.Range("E23").Value = '  +1.96%  '
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("E25").Value = '  +0.96%  '
$ws.Range("E26").Value = '  +8.38%  '
$ws.Range("E27").Value = '  -0.86%  '
$ws.Range("E28").Value = '  +1.27%  '
$ws.Range("E29").Value = '  +0.64%  '
$ws.Range("E30").Value = '  +2.27%  '
$ws.Range("E31").Value = '  +1.20%  '
$ws.Range("E32").Value = '  +0.67%  '
$ws.Range("E33").Value = '  +8.46%  '
$ws.Range("E34").Value = '  +0.80%  '
$ws.Range("E35").Value = '  +5.99%  '
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("E37").Value = '  +3.84%  '
$ws.Range("E38").Value = '  +3.08%  '
$ws.Range("E39").Value = '  +2.47%  '
$ws.Range("E40").Value = '  -0.41%  '
$ws.Range("E41").Value = '  +0.46%  '
$ws.Range("E42").Value = '  +3.05%  '
$ws.Range("E43").Value = '  +0.55%  '
$ws.Range("E44").Value = '  +1.99%  '
$ws.Range("E45").Value = '  +1.25%  '
$ws.Range("E46").Value = '  +1.26%  '
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("E48").Value = '  +0.71%  '
$ws.Range("E49").Value = '  +3.08%  '
$ws.Range("E50").Value = '  +1.79%  '
$ws.Range("E51").Value = '  +0.55%  '
